# Optuna Attempt (go back with original)
# Updates forecast metrics (Inventory Coverage / Seasonality Index) on the
# "Forecast Comparison" sheet and the rolling forecast totals on the
# "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Column H = Inventory Coverage, Column L = Seasonality Index

$inventoryCoverage = @{
    2 = 3.85
    3 = 2.85
    4 = 1.85
    5 = 0.85
}

foreach ($row in $inventoryCoverage.Keys) {
    $wsForecast.Range("H$row").Value = $inventoryCoverage[$row]
}

$seasonalityIndex = @{
    2  = 0.86
    3  = 1.14
    4  = 0.88
    5  = 0.99
    6  = 1.07
    7  = 0.83
    8  = 1.17
    9  = 0.88
    10 = 1.07
    11 = 0.9399999999999999
    12 = 0.91
    13 = 1.14
    14 = 1.04
    15 = 1.15
    16 = 0.88
    17 = 0.87
}

foreach ($row in $seasonalityIndex.Keys) {
    $wsForecast.Range("L$row").Value = $seasonalityIndex[$row]
}

# --- Summary sheet ---------------------------------------------------------
# These cells store numeric-looking values as TEXT, so we force the cell to
# text format before assigning and then drop back to the Normal style so we
# don't leave a stray number-format override behind.

function Set-TextCellValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextCellValue $wsSummary.Range("B9")  "20"
Set-TextCellValue $wsSummary.Range("B10") "10"
Set-TextCellValue $wsSummary.Range("B11") "5"
